# Apply the updated cryptos list values (Price column D, Volume(1h) column E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.040.35"
$ws.Range("E2").Value = "  +1.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.265.67"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.93"
$ws.Range("E6").Value = "  +2.53%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +1.13%  "
$ws.Range("E9").Value = "  -2.11%  "
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("E11").Value = "  -2.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.843.75"
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.51"
$ws.Range("E14").Value = "  -1.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "68.039.84"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.181.17"
$ws.Range("E17").Value = "  -6.17%  "
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "415.98"
$ws.Range("E20").Value = "  +7.41%  "
$ws.Range("E21").Value = "  -3.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.41"
$ws.Range("E23").Value = "  +0.95%  "
$ws.Range("E24").Value = "  -1.34%  "
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  -0.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.66"
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.46"
$ws.Range("E31").Value = "  -2.61%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("E33").Value = "  -2.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.24"
$ws.Range("E34").Value = "  -1.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "162.51"
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.45"
$ws.Range("E36").Value = "  -2.42%  "
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.05"
$ws.Range("E38").Value = "  +1.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.797"
$ws.Range("E39").Value = "  -1.81%  "
$ws.Range("E40").Value = "  -2.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.36"
$ws.Range("E41").Value = "  -1.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.637.38"
$ws.Range("E42").Value = "  +1.23%  "
$ws.Range("E43").Value = "  -1.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0676"
$ws.Range("E44").Value = "  -0.86%  "
$ws.Range("E45").Value = "  -2.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "337.04"
$ws.Range("E46").Value = "  -0.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.31"
$ws.Range("E47").Value = "  -1.28%  "
$ws.Range("E48").Value = "  -1.88%  "
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.977"
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("E51").Value = "  -1.44%  "
